# Generate Report for Handoff
# Adds two new localization entries (2d097034-3bf1-43c8-b50e-fad5b6f4ebc4 and
# d945e330-49de-42ad-bed2-f24a7e3cf67e) as new rows 6 & 7 on all three sheets:
# "Overview", "zh-cn", "de-de".

$wb = $excel.ActiveWorkbook

$entries = @(
    @{
        Guid        = "2d097034-3bf1-43c8-b50e-fad5b6f4ebc4"
        MdHash      = "2d097034-3bf1-43c8-b50e-fad5b6f4ebc4"
        XlfHashZh   = "58f8e63d6d8475ece6ef72ca76bde4242bdb2875"
        XlfHashDe   = "58f8e63d6d8475ece6ef72ca76bde4242bdb2875"
        HandoffDate = "2016-29-14 06:29:41"
        ZhXlfDate   = "2016-03-14 06:29:38"
        DeXlfDate   = "2016-03-14 06:29:41"
    },
    @{
        Guid        = "d945e330-49de-42ad-bed2-f24a7e3cf67e"
        MdHash      = "d945e330-49de-42ad-bed2-f24a7e3cf67e"
        XlfHashZh   = "11d67288cfd3aa369ac443656864310b4b643ad9"
        XlfHashDe   = "11d67288cfd3aa369ac443656864310b4b643ad9"
        HandoffDate = "2016-29-14 06:29:41"
        ZhXlfDate   = "2016-03-14 06:29:41"
        DeXlfDate   = "2016-03-14 06:29:41"
    }
)

$status = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewStartRow = 6

for ($i = 0; $i -lt $entries.Count; $i++) {
    $e = $entries[$i]
    $r = $overviewStartRow + $i
    $mdName = $e.Guid + ".md"

    $wsOverview.Range("B$r").Value = $status
    $wsOverview.Range("C$r").Value = $status
    $wsOverview.Range("D$r").Value = $e.HandoffDate

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/" + $e.MdHash + "/e2e/" + $mdName
    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$r"), $mdUrl, "", "", $mdName)
}

# ---------------------------------------------------------------------
# Language sheets "zh-cn" / "de-de":
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; Lang = "zh-cn"; HashField = "XlfHashZh"; DateField = "ZhXlfDate"; OrgRepo = "oltest.zh-cn" },
    @{ Name = "de-de"; Lang = "de-de"; HashField = "XlfHashDe"; DateField = "DeXlfDate"; OrgRepo = "oltest.de-de" }
)

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $startRow = 6

    for ($i = 0; $i -lt $entries.Count; $i++) {
        $e = $entries[$i]
        $r = $startRow + $i
        $mdName = $e.Guid + ".md"
        $xlfHash = $e[$lang.HashField]
        $xlfDate = $e[$lang.DateField]
        $xlfName = $e.Guid + "." + $xlfHash + "." + $lang.Lang + ".xlf"

        $ws.Range("C$r").Value = $status
        $ws.Range("E$r").Value = $xlfDate
        $ws.Range("H$r").Value = "0001-01-01 00:00:00"
        $ws.Range("I$r").Value = "Include"

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/" + $e.MdHash + "/e2e/" + $mdName
        $ws.Hyperlinks.Add($ws.Range("A$r"), $mdUrl, "", "", $mdName)
        $ws.Hyperlinks.Add($ws.Range("B$r"), $mdUrl, "", "", ".md")

        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/" + $xlfHash + "/ol-handoff/OpenLocalizationTestOrg/" + $lang.OrgRepo + "/ci/ht/" + $xlfName
        $ws.Hyperlinks.Add($ws.Range("D$r"), $xlfUrl, "", "", $xlfName)
    }
}

Write-Host "Generated handoff report rows for 2d097034-3bf1-43c8-b50e-fad5b6f4ebc4 and d945e330-49de-42ad-bed2-f24a7e3cf67e"
